# Insert a new weekly data row at row 56 (pushing existing rows 56..86 down to 57..87)
# and populate it with this week's Jengibre price data for Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 56; this shifts rows 56-86 to 57-87
# and extends the sheet dimension to A1:R87.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Range("A56").Value = 8
$ws.Range("B56").Value = "Terminal La Palmera de La Serena"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value2 = 44917
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100114007
$ws.Range("G56").Value = "Jengibre"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 540
$ws.Range("K56").Value = 14000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = 14500
$ws.Range("N56").Value = "$/caja 13 kilos"
$ws.Range("O56").Value = "Perú"
$ws.Range("P56").Value = 1115
$ws.Range("Q56").Value = 13
$ws.Range("R56").Value = "Hortaliza"
